# Attendance taking feature added to the system
# - Update the date header in D1
# - Replace the student roster (rows 3-13) with the new attendance list
# - Remove the old trailing rows (14-17)
# - Mark attendance (D column = 1) for every present student; row 9 (roll 80,
#   Shubham Pandey) is left blank to represent an absence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe out the old roster (rows 3-17) first so stale names/rolls don't linger
# in the sheet beyond the new, smaller attendance list.
$ws.Range("A3:D17").Clear()

# Update the date column header
$ws.Range("D1").Value = "28_12_18"

# New roster data: Roll Number, First Name, Last Name
$data = @(
    @("104", "Yash", "Atre"),
    @("22", "Devender", "Singh"),
    @("34", "krutik", "pathak"),
    @("44", "Manas", "Jain"),
    @("7", "Akshat", "Gupta"),
    @("73", "Sakina", "Saifee"),
    @("80", "Shubham", "Pandey"),
    @("85", "Shweta", "Solanki"),
    @("9", "Aman", "Bhawsar"),
    @("95", "syed", "mustafa"),
    @("55", "palak", "agrawal")
)

$startRow = 3
$endRow = $startRow + $data.Count - 1

# Roll numbers are stored as text (e.g. "104", "22"), so format column A as
# text first, otherwise Excel would coerce the numeric-looking strings to
# actual numbers on assignment.
$rollRange = "A" + $startRow + ":A" + $endRow
$ws.Range($rollRange).NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    if ($row -eq 9) {
        # absent student, no attendance mark
        $ws.Cells.Item($row, 4).Value = $null
    } else {
        $ws.Cells.Item($row, 4).Value = 1
    }
}
